$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F, shifting EMPLOYEE_ID..TYPE_OF_SALARY (old F:N)
# one column to the right (new G:O).
$ws.Columns("F:F").Insert()

# New column F header + blank data cells (pandas "Unnamed: 0" index column).
$ws.Range("F1").Value = "Unnamed: 0"
$ws.Range("F2").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("F5").Value = ""

# Column H (old G / MANAGER_ID) gets refreshed manager-id values.
$ws.Range("H2").Value = "O254"
$ws.Range("H3").Value = "O281"
$ws.Range("H4").Value = "O254"
$ws.Range("H5").Value = "O254"

# Row 3 (SAKSHI KHANNA): status/process/department updated.
$ws.Range("L3").Value = "INACTIVE"
$ws.Range("M3").Value = "HERO"
$ws.Range("N3").Value = "RECOVERY"

# Row 5 (VINITA KUMARI): process/department updated.
$ws.Range("M5").Value = "IDFC"
$ws.Range("N5").Value = "TW"
